# "Level" (column E) for the Masters-level rows (38-52) was stored as the
# literal text "M". Replace it with the actual numeric level: 4 for rows
# 38-44, 5 for rows 45-52.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E38:E44").Value = 4
$ws.Range("E45:E52").Value = 5

# Reflect where the author ended up looking in the sheet (best effort -
# scroll position isn't always persisted, but set it anyway).
try {
    $excel.ActiveWindow.ScrollRow = 31
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

$ws.Range("G45:G52").Select()
